$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reflect the resized Excel window (best effort - some hosts expose this
# as read-only, so failures here are harmless and ignored).
try {
    $win = $wb.Windows.Item(1)
    $win.Width = 28800
    $win.Height = 13170
} catch {}

# Add 3 new groups of data (columns E, F, G)
# Shared-string table must end up in order: cxq, hyy, hzj (indices 4,5,6)
$ws.Range("E1").Value = "cxq6hz_20170224_144343_ASIC_EEG"
$ws.Range("G1").Value = "hyy-调节6Hz_20170306_110203_ASIC_EEG"
$ws.Range("F1").Value = "hzj-调节6Hz_20170220_113105_ASIC_EEG"

$ws.Range("E2").Value = 0.96142433234421365
$ws.Range("F2").Value = 0.94594594594594594
$ws.Range("G2").Value = 0.91919191919191923

$ws.Range("E3").Value = 0.92128279883381925
$ws.Range("F3").Value = 0.9
$ws.Range("G3").Value = 0.94197952218430037

$ws.Columns.Item(6).Select()
